# Apply the "1.2 -> 1.3" edit to the RF008 workbook.
# Summary of change:
#   - TC3's old steps 6 ("verifica que os campos...limpos") and the old
#     step 7 ("Apontamentos"/"apresenta os Apontamentos") are removed
#     (rows 39-42 deleted), shifting the TC4 block (rows 45-57) up by 4
#     rows to become rows 41-53.
#   - TC3's final remaining step (now row 38, step 5) becomes "clica na
#     opcao 'Salvar'" with expected result "erro ao tentar salvar".
#   - TC4's final step (now row 53, step 8) expected result becomes
#     "erro ao tentar editar" (swapped with TC3's old message).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Delete the 4 rows that held TC3's old steps 6 and 7 (rows 39-42).
#    This shifts everything below (the TC4 block, previously rows 45-57)
#    up by 4 rows, so it becomes rows 41-53.
$ws.Range("A39:F42").EntireRow.Delete() | Out-Null

# 2. Update TC3's now-final step (row 38, step 5) to be the "Salvar" step
#    with the "erro ao tentar salvar" expected result.
$ws.Range("B38").Value = "Avaliador de Pessoas clica na opcao 'Salvar'"
$ws.Range("D38").Value = "SYSTEM exibe uma mensagem de erro ao tentar salvar, informando o campo ou a validacao que falhou"

# 3. Update TC4's final step (row 53, step 8) expected result to
#    "erro ao tentar editar".
$ws.Range("D53").Value = "SYSTEM exibe uma mensagem de erro ao tentar editar, informando o campo ou a validacao que falhou"
